# Fix: comments should be associated with (and reference) the sheet they
# belong to. Sheet1's comments are re-labelled to make their owning sheet
# explicit, a new "Different Name" sheet is added with its own data +
# comments, and a B2 comment is added to Sheet1.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Sheet1: update existing comment text, add a new one at B2 ---------
$ws1.Range("A2").Comment.Text("A2Sheet1: This is a random sample comment") | Out-Null
$ws1.Range("B2").AddComment("B2Sheet1: Empty Value") | Out-Null
$ws1.Range("A3").Comment.Text("A3Sheet1: Another comment") | Out-Null
$ws1.Range("A5").Comment.Text("A5Sheet1: Comment with different font and size") | Out-Null

# --- Add the new worksheet right after Sheet1 ---------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Different Name"

# --- Populate the new sheet's data --------------------------------------
$ws2.Range("A1").Value = "1 Value"
$ws2.Range("B2").Value = "2 Value"
$ws2.Range("A3").Value = "Value 3"
$ws2.Range("C3").Value = "3 Value"
$ws2.Range("A4").Value = "Value 4"
$ws2.Range("D4").Value = "4 Value"
$ws2.Range("A5").Value = "Value 5"
$ws2.Range("E5").Value = "5 Value"

# --- Add comments on the new sheet, each naming its own cell/sheet -----
$ws2.Range("A2").AddComment("A2DN: Empty Value") | Out-Null
$ws2.Range("B2").AddComment("B2DN: Random Comment") | Out-Null
$ws2.Range("A3").AddComment("A3DN: Another comment") | Out-Null
$ws2.Range("C3").AddComment("C3DN: C3") | Out-Null
$ws2.Range("A5").AddComment("A5DN: Comment with different font and size") | Out-Null

# Leave Sheet1 as the active / selected sheet.
$ws1.Activate()
